$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the header text values in A1:N1, keeping their existing style/formatting.
$ws.Range("A1:N1").ClearContents()

# Update the active selection to N4 (matches the saved view state in the diff).
$ws.Range("N4").Select()
